# Insert a new data row at row 30 (pushing existing rows 30-63 down to 31-64)
# and populate it with the new weekly price record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(30).Insert()

$ws.Range("A30").Value = 9
$ws.Range("B30").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C30").Value = "Metropolitana"
$ws.Range("D30").Value = 45036
$ws.Range("E30").Value = 13
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100104
$ws.Range("H30").Value = "Frutos de pepita"
$ws.Range("I30").Value = 100104003
$ws.Range("J30").Value = "Membrillo"
$ws.Range("K30").Value = "Champion"
$ws.Range("L30").Value = "Primera"
$ws.Range("M30").Value = 18
$ws.Range("N30").Value = 260000
$ws.Range("O30").Value = 270000
$ws.Range("P30").Value = 265556
$ws.Range("Q30").Value = "`$/bins (450 kilos)"
$ws.Range("R30").Value = "Región de O'Higgins"
$ws.Range("S30").Value = 590
$ws.Range("T30").Value = 450
